# The worksheet is protected (legacy password hash) so we must unprotect it
# before any cell can be modified, then re-protect it afterwards to leave the
# sheet in the same protected state it started in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the "as of" date in the confidentiality / disclosure note (A18),
# from 2021-06-10 to 2021-06-14.
$null = $ws.Cells.Replace("2021-06-10", "2021-06-14")

# That note wraps onto two lines, and re-writing it makes Excel stamp an
# explicit custom row height on row 18. Auto-fit the row back so it matches
# the original (no explicit height) state.
$ws.Rows("18").AutoFit()

# Update the Weight (D) and Percent Change (E) columns for rows 2-15.
$ws.Range("D2").Value = 0.05778225121240245
$ws.Range("E2").Value = 0.001538322323771713

$ws.Range("D3").Value = 0.02086639390334767
$ws.Range("E3").Value = 0.001503476790077007

$ws.Range("D4").Value = 0.02858001856184886
$ws.Range("E4").Value = 0.004024881083058984

$ws.Range("D5").Value = 0.03019960523391851
$ws.Range("E5").Value = 0.002207505518763808

$ws.Range("D6").Value = 0.03083619821963111
$ws.Range("E6").Value = -0.009325985587113172

$ws.Range("D7").Value = 0.01866642266114168
$ws.Range("E7").Value = -0.01249999999999996

$ws.Range("D8").Value = 0.01030770839596868
$ws.Range("E8").Value = 0.01198402130492671

$ws.Range("D9").Value = 0.01043515771035673
$ws.Range("E9").Value = -0.0005636978579481866

$ws.Range("D10").Value = 0.0694422295133397
$ws.Range("E10").Value = -0.005364806866952843

$ws.Range("D11").Value = 0.0695539927582646
$ws.Range("E11").Value = -0.005356186395286566

$ws.Range("D12").Value = 0.1488189696866708
$ws.Range("E12").Value = -0.007659335254022959

$ws.Range("D13").Value = 0.3929177396373904
$ws.Range("E13").Value = -0.003113917481186834

$ws.Range("D14").Value = 0.1115933125057189
$ws.Range("E14").Value = -0.0008199601733630546

$ws.Range("E15").Value = -0.003301264035764251

# Re-protect the sheet to restore its original protected state.
$ws.Protect()
